{"js": "// Commit: \"update to verse 2\"\n//   \"With fury whipped with schemes so grim\"  -> \"With fury whipped and schemes so grim\"\n//   \"They scream their bile and hate at him.\" -> \"To scream their bile and hate at him.\"\n\nconst body = context.document.body;\n\n// --- Edit 1: the whole word \"with\" -> \"and\" (\"...whipped with schemes...\" -> \"...whipped and schemes...\") ---\nconst withResults = body.search(\"with\", { matchCase: true, matchWholeWord: true });\nwithResults.load(\"items\");\nawait context.sync();\n\nfor (let i = withResults.items.length - 1; i >= 0; i--) {\n  withResults.items[i].insertText(\"and\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- Edit 2: \"They \" -> \"To \" (\"They scream...\" -> \"To scream...\") ---\nconst theyResults = body.search(\"They \", { matchCase: true });\ntheyResults.load(\"items\");\nawait context.sync();\n\nfor (let i = theyResults.items.length - 1; i >= 0; i--) {\n  theyResults.items[i].insertText(\"To \", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Commit: \"update to verse 2\"\n#   \"With fury whipped with schemes so grim\"  -> \"With fury whipped and schemes so grim\"\n#   \"They scream their bile and hate at him.\" -> \"To scream their bile and hate at him.\"\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: the whole word \"with\" -> \"and\" ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\"with\", $true, $true, $false, $false, $false, $true, 1, $false, \"and\", 1)\n\n# --- Edit 2: \"They \" -> \"To \" ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\"They \", $true, $false, $false, $false, $false, $true, 1, $false, \"To \", 1)\n"}
